# Updates cryptos list values (Price and Volume(1h) columns) to match the
# latest scraped snapshot. Cells keep their original "text" storage so
# values such as "511.00" or "67.095.47" are preserved verbatim instead of
# being reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.095.47"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.485.60"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.19%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "583.37"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "171.25"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "2.485.34"
$ws.Range("E9").Value = "  +2.72%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.138"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +4.42%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "2.912.03"
$ws.Range("E14").Value = "  +1.40%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "25.34"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "67.133.63"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "2.481.51"
$ws.Range("E18").Value = "  +1.96%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "10.99"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -2.48%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.43"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -2.25%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "348.47"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -0.81%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "68.40"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +0.68%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "9.31"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +2.35%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "0.0₃0907"
$ws.Range("E30").Value = "  +0.24%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "511.00"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -0.32%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "160.86"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.117"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("E38").Value = "  +1.01%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "18.24"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("E44").Value = "  +1.33%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "2.37"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "38.82"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "142.88"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("E48").Value = "  -1.17%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.515"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +0.23%  "
